# Added View Profile functionality
#
# Insert a new empty paragraph (Normal style) immediately after the
# paragraph containing the "ADD COLUMN `rating` ..." SQL statement,
# matching the author's diff which added one blank paragraph there.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*ADD COLUMN*rating*DOUBLE NOT NULL AFTER*longitude*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()
    Write-Output "Inserted new paragraph after rating column line."
} else {
    Write-Output "Target paragraph not found."
}
